$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# K2: Signal/Method header for second repetition's last column changed
# from "scaledILD_0.6" to "scaledILD_1"
$ws.Range("K2").Value = "scaledILD_1"

# Data entry: subjects rated "90" for a number of cells that were
# previously "0" (unprocessed placeholders)
$ws.Range("H3").Value = 90
$ws.Range("H4").Value = 90
$ws.Range("D5").Value = 90
$ws.Range("E6").Value = 90
$ws.Range("F7").Value = 90
$ws.Range("H9").Value = 90
$ws.Range("E10").Value = 90
$ws.Range("C11").Value = 90
$ws.Range("H12").Value = 90
$ws.Range("G13").Value = 90
$ws.Range("E16").Value = 90
$ws.Range("G17").Value = 90
$ws.Range("J18").Value = 90
$ws.Range("H19").Value = 90
$ws.Range("E20").Value = 90
$ws.Range("H22").Value = 90
$ws.Range("G23").Value = 90
$ws.Range("F24").Value = 90
$ws.Range("E25").Value = 90
$ws.Range("G26").Value = 90
$ws.Range("H26").Value = 90

# The comment rows (7, 13, 20, 26) shrink from 135pt to 112pt tall
$ws.Rows.Item(7).RowHeight = 112
$ws.Rows.Item(13).RowHeight = 112
$ws.Rows.Item(20).RowHeight = 112
$ws.Rows.Item(26).RowHeight = 112

# Selection / view moved while editing
$ws.Range("M5").Select()
